$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo: remove stray trailing parenthesis in E9
$ws.Range("E9").Value = "Descriptive Statistics."

# Update the active cell / selection to match the saved view state
$ws.Range("E14").Select()
